$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove comment 0 ("Pattern" heading comment) and the " und HIWEISE?" text
#    that was appended next to "Pattern".
# ---------------------------------------------------------------------------
$d.Comments.Item(1).Delete()

$rPattern = $d.Paragraphs(6).Range
$rPattern.Find.Execute(" und HIWEISE?", $false, $false, $false, $false, $false, `
    $true, 1, $false, "", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Re-add the "_GoBack" bookmark right after the word "Pattern" (this also
#    implicitly removes the old "_GoBack" bookmark that used to sit near
#    "Titel tragen koennen", since bookmark names are unique per document).
#    A temporary trailing space is used as a workaround so the bookmark is
#    not anchored exactly at the paragraph end (a position that otherwise
#    gets mis-serialized to the next paragraph).
# ---------------------------------------------------------------------------
$pPattern = $d.Paragraphs(6).Range
$pPattern.Collapse(0)
$pPattern.MoveEnd(1, -1) | Out-Null
$pPattern.InsertAfter(" ")

$pPattern2 = $d.Paragraphs(6).Range
$bmPos = $pPattern2.End - 2
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
$d.Range($bmPos, $bmPos + 1).Delete()

# ---------------------------------------------------------------------------
# 3) Move the "Auf die Labels der Relationen..." paragraph from right after
#    the Koordinatorpattern paragraph to right after "Hinweis" paragraph
#    (i.e. immediately before the "Pattern" heading).
# ---------------------------------------------------------------------------
$rHinweisBody = $d.Paragraphs(5).Range
$rHinweisBody.Find.Execute("nicht modelliert.", $false, $false, $false, $false, $false, `
    $true, 1, $false, ("nicht modelliert.^pAuf die Labels der Relationen wurde aus " + `
    "Übersichtszwecken verzichtet. Außerdem werden die Attribute die einer Rolle " + `
    "zugewiesen sind immer auch als Attribut vermerkt, d.h. die Rollen ersetzten keine Attribute"), `
    2) | Out-Null

# Delete the now-duplicated original "Auf die Labels..." paragraph (it sits
# right after the Koordinatorpattern paragraph, i.e. paragraph 9 once the
# new paragraph 6 above has been inserted).
$d.Paragraphs(9).Range.Delete()

# ---------------------------------------------------------------------------
# 4) Split the Rollen-/Koordinatorpattern paragraph into three paragraphs and
#    flesh out the Rollenpattern explanation.
# ---------------------------------------------------------------------------
$rRollen = $d.Paragraphs(8).Range
$rRollen.Find.Execute("verwendet. Da einige Klassen mehr als eine Referenz auf eine Klasse haben, bietet sich hier dieses Pattern an. Das Koordinatorpattern", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    ("verwendet. ^pDas Rollenpattern bietet sich deswegen an, da einige Klassen mehr als eine " + `
    "Referenz auf eine Klasse haben. Daher agieren sie hierbei als Rollen mit entsprechenden Rollennamen.^p" + `
    "Das Koordinatorpattern"), 2) | Out-Null

# ---------------------------------------------------------------------------
# 5) Remove comment 1 from the "Person" paragraph, merging the commented text
#    back into the surrounding sentence (same wording, single run).
# ---------------------------------------------------------------------------
$d.Comments.Item(1).Delete()

